$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("configuration")

# --- Add new column H (report_receiver) ---
# Copy formats from similarly-styled existing cells before filling in the new values.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122) # xlPasteFormats
$ws.Range("C2").Copy()
$ws.Range("H2").PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = 0

$ws.Range("H1").Value = "report_receiver"

# --- Update existing data values ---
$ws.Range("G2").Value = "Both"

# --- Rename the first header column ---
$ws.Range("A1").Value = "project_name"

$ws.Range("E2").Value = "yes"

$ws.Range("H2").Value = "unidevgo.qa1@gmail.com"

$ws.Columns.Item(8).ColumnWidth = 23.6

# --- Update the data validation range for the yes/no list (now a contiguous E2:F2 range) ---
$ws.Range("E2:F2").Validation.Delete()
$ws.Range("E2:F2").Validation.Add(3, 1, 1, """no, yes""")
$ws.Range("E2:F2").Validation.IgnoreBlank = $false
$ws.Range("E2:F2").Validation.ShowInput = $true
$ws.Range("E2:F2").Validation.ShowError = $true

# --- Update the selected cell ---
$ws.Range("G2").Select()
